$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 updates
$ws.Range("G5").Value = 2.15
$ws.Range("I5").Value = 3.3
$ws.Range("J5").Value = 3.1
$ws.Range("L5").Value = 4.33
$ws.Range("O5").Value = 1.5
$ws.Range("P5").Value = 2.5
$ws.Range("Q5").Value = 1.95
$ws.Range("R5").Value = 1.9
$ws.Range("S5").Value = 2.5
$ws.Range("T5").Value = 1.5
$ws.Range("U5").Value = 5
$ws.Range("V5").Value = 1.17
$ws.Range("W5").Value = 1.57
$ws.Range("X5").Value = 2.25
$ws.Range("Y5").Value = 2.2
$ws.Range("Z5").Value = 1.62
$ws.Range("AA5").Value = 6
$ws.Range("AB5").Value = 9
$ws.Range("AD5").Value = 21
$ws.Range("AG5").Value = 6.5
$ws.Range("AI5").Value = 19
$ws.Range("AJ5").Value = 67
$ws.Range("AK5").Value = 8
$ws.Range("AL5").Value = 15
$ws.Range("AM5").Value = 13
$ws.Range("AO5").Value = 34
$ws.Range("AP5").Value = 41
$ws.Range("AQ5").Value = 1250
$ws.Range("AR5").Value = 4.2
$ws.Range("AS5").Value = 1.22

# Row 6 updates
$ws.Range("G6").Value = 1.75
$ws.Range("J6").Value = 2.4
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 10
$ws.Range("AC6").Value = 8.5
$ws.Range("AH6").Value = 6.5

# Row 8 updates
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 7
